$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rule name in B11 ("R40") is replaced with the text "1". A plain
# Range.Value = "1" would be auto-coerced to a number by the object model,
# which would change both the cell's stored type and its style (Excel
# needs a distinct number-format/quote-prefix flavoured style to hold a
# numeric-looking string). To land the new content as genuine text while
# leaving B11's existing style untouched, compute "1" as a formula result
# (always a string/text value) in a scratch cell, then copy/paste just the
# value into B11 - this keeps B11's own formatting exactly as it was.
$helper = $ws.Range("ZZ1")
$helper.Formula = '="1"'
$helper.Copy()

$target = $ws.Range("B11")
$target.PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = $false
$helper.Clear()
